$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3 values per diff
$ws.Range("G3").Value = 3
$ws.Range("H3").Value = 3.4
$ws.Range("I3").Value = 2.3
$ws.Range("J3").Value = 3.6
$ws.Range("L3").Value = 3
$ws.Range("X3").Value = 15
$ws.Range("AA3").Value = 23
$ws.Range("AB3").Value = 34
$ws.Range("AI3").Value = 11
$ws.Range("AJ3").Value = 9.5
$ws.Range("AK3").Value = 21
$ws.Range("AL3").Value = 19
$ws.Range("AN3").Value = 5
$ws.Range("AO3").Value = 17
$ws.Range("AP3").Value = 26
$ws.Range("AW3").Value = 4.33

# Delete row 4 entirely (shrinks used range so dimension becomes A1:BB3)
$ws.Rows("4:4").Delete()
